$wb = $excel.ActiveWorkbook

# --- Update "Last Updated" timestamp on Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "04 Nov 2025, 09:20 AM"

# --- Refresh Stock List sheet data (rows 2-76) ---
$ws = $wb.Worksheets.Item("Stock List")

$bcde = New-Object 'object[,]' 75,4
$bcde[0,0] = "CANHLIFE"
$bcde[0,1] = "CANHLIFE"
$bcde[0,2] = 119
$bcde[0,3] = 1.0873
$bcde[1,0] = "RUBICON"
$bcde[1,1] = "RUBICON"
$bcde[1,2] = 662
$bcde[1,3] = 1.2852
$bcde[2,0] = "CRAMC"
$bcde[2,1] = "CRAMC"
$bcde[2,2] = 309.8
$bcde[2,3] = -0.0645
$bcde[3,0] = "LGEINDIA"
$bcde[3,1] = "LGEINDIA"
$bcde[3,2] = 1645.1
$bcde[3,3] = -0.2365
$bcde[4,0] = "TATACAP"
$bcde[4,1] = "TATACAP"
$bcde[4,2] = 326.4
$bcde[4,3] = -0.7299
$bcde[5,0] = "WEWORK"
$bcde[5,1] = "WEWORK"
$bcde[5,2] = 643
$bcde[5,3] = -0.7256
$bcde[6,0] = "ADVANCE"
$bcde[6,1] = "ADVANCE"
$bcde[6,2] = 137
$bcde[6,3] = -0.204
$bcde[7,0] = "OMFREIGHT"
$bcde[7,1] = "OMFREIGHT"
$bcde[7,2] = 89.05
$bcde[7,3] = -0.4249
$bcde[8,0] = "GLOTTIS"
$bcde[8,1] = "GLOTTIS"
$bcde[8,2] = 74
$bcde[8,3] = 0.8587
$bcde[9,0] = "FABTECH"
$bcde[9,1] = "FABTECH"
$bcde[9,2] = 239.85
$bcde[9,3] = 1.3736
$bcde[10,0] = "PACEDIGITK"
$bcde[10,1] = "PACEDIGITK"
$bcde[10,2] = 218.9
$bcde[10,3] = 0.1556
$bcde[11,0] = "JAINREC"
$bcde[11,1] = "JAINREC"
$bcde[11,2] = 372.65
$bcde[11,3] = -0.0134
$bcde[12,0] = "EPACKPEB"
$bcde[12,1] = "EPACKPEB"
$bcde[12,2] = 291
$bcde[12,3] = -1.5562
$bcde[13,0] = "BMWVENTLTD"
$bcde[13,1] = "BMWVENTLTD"
$bcde[13,2] = 66.81
$bcde[13,3] = -3.5235
$bcde[14,0] = "STYL"
$bcde[14,1] = "STYL"
$bcde[14,2] = 375
$bcde[14,3] = -0.1464
$bcde[15,0] = "JARO"
$bcde[15,1] = "JARO"
$bcde[15,2] = 630
$bcde[15,3] = -0.1347
$bcde[16,0] = "SOLARWORLD"
$bcde[16,1] = "SOLARWORLD"
$bcde[16,2] = 312.65
$bcde[16,3] = 0.5144
$bcde[17,0] = "ARSSBL"
$bcde[17,1] = "ARSSBL"
$bcde[17,2] = 515
$bcde[17,3] = 0.3801
$bcde[18,0] = "GANESHCP"
$bcde[18,1] = "GANESHCP"
$bcde[18,2] = 282.3
$bcde[18,3] = 0
$bcde[19,0] = "ATLANTAELE"
$bcde[19,1] = "ATLANTAELE"
$bcde[19,2] = 1019
$bcde[19,3] = -0.1812
$bcde[20,0] = "GKENERGY"
$bcde[20,1] = "GKENERGY"
$bcde[20,2] = 216
$bcde[20,3] = 0.2041
$bcde[21,0] = "SAATVIKGL"
$bcde[21,1] = "SAATVIKGL"
$bcde[21,2] = 535.2
$bcde[21,3] = 0
$bcde[22,0] = "IVALUE"
$bcde[22,1] = "IVALUE"
$bcde[22,2] = 280.2
$bcde[22,3] = -0.779
$bcde[23,0] = "VMSTMT"
$bcde[23,1] = "VMSTMT"
$bcde[23,2] = 70.51
$bcde[23,3] = -0.2264
$bcde[24,0] = "EUROPRATIK"
$bcde[24,1] = "EUROPRATIK"
$bcde[24,2] = 319.15
$bcde[24,3] = 0
$bcde[25,0] = "SHRINGARMS"
$bcde[25,1] = "SHRINGARMS"
$bcde[25,2] = 233.68
$bcde[25,3] = 0.62
$bcde[26,0] = "DEVX"
$bcde[26,1] = "DEVX"
$bcde[26,2] = 44.95
$bcde[26,3] = 0.5593
$bcde[27,0] = "URBANCO"
$bcde[27,1] = "URBANCO"
$bcde[27,2] = 151.6
$bcde[27,3] = -0.2697
$bcde[28,0] = "AMANTA"
$bcde[28,1] = "AMANTA"
$bcde[28,2] = 120.01
$bcde[28,3] = -0.6704
$bcde[29,0] = "CPEDU"
$bcde[29,1] = "CPEDU"
$bcde[29,2] = 313
$bcde[29,3] = 0.9189
$bcde[30,0] = "AHCL"
$bcde[30,1] = "AHCL"
$bcde[30,2] = 132.85
$bcde[30,3] = -1.5853
$bcde[31,0] = "STLNETWORK"
$bcde[31,1] = "STLNETWORK"
$bcde[31,2] = 26.69
$bcde[31,3] = -0.0375
$bcde[32,0] = "VIKRAN"
$bcde[32,1] = "VIKRAN"
$bcde[32,2] = 100.4
$bcde[32,3] = 0.571
$bcde[33,0] = "MEIL"
$bcde[33,1] = "MEIL"
$bcde[33,2] = 468.95
$bcde[33,3] = 0.9472
$bcde[34,0] = "SHREEJISPG"
$bcde[34,1] = "SHREEJISPG"
$bcde[34,2] = 273.1
$bcde[34,3] = 0.3306
$bcde[35,0] = "GEMAROMA"
$bcde[35,1] = "GEMAROMA"
$bcde[35,2] = 221.47
$bcde[35,3] = 0.0045
$bcde[36,0] = "PATELRMART"
$bcde[36,1] = "PATELRMART"
$bcde[36,2] = 221.67
$bcde[36,3] = 0
$bcde[37,0] = "VIKRAMSOLR"
$bcde[37,1] = "VIKRAMSOLR"
$bcde[37,2] = 325.6
$bcde[37,3] = -0.489
$bcde[38,0] = "REGAAL"
$bcde[38,1] = "REGAAL"
$bcde[38,2] = 90.36
$bcde[38,3] = 0.5005
$bcde[39,0] = "BLUESTONE"
$bcde[39,1] = "BLUESTONE"
$bcde[39,2] = 714.5
$bcde[39,3] = 0.4852
$bcde[40,0] = "ALLTIME"
$bcde[40,1] = "ALLTIME"
$bcde[40,2] = 300
$bcde[40,3] = -0.2494
$bcde[41,0] = "JSWCEMENT"
$bcde[41,1] = "JSWCEMENT"
$bcde[41,2] = 135.5
$bcde[41,3] = -0.0959
$bcde[42,0] = "HILINFRA"
$bcde[42,1] = "HILINFRA"
$bcde[42,2] = 77.96
$bcde[42,3] = 0.5417
$bcde[43,0] = "LOTUSDEV"
$bcde[43,1] = "LOTUSDEV"
$bcde[43,2] = 177.23
$bcde[43,3] = 0.0339
$bcde[44,0] = "MBEL"
$bcde[44,1] = "MBEL"
$bcde[44,2] = 453.7
$bcde[44,3] = 0
$bcde[45,0] = "LAXMIINDIA"
$bcde[45,1] = "LAXMIINDIA"
$bcde[45,2] = 148.8
$bcde[45,3] = 0.9635
$bcde[46,0] = "CPPLUS"
$bcde[46,1] = "CPPLUS"
$bcde[46,2] = 1329.8
$bcde[46,3] = 0.3168
$bcde[47,0] = "SHANTIGOLD"
$bcde[47,1] = "SHANTIGOLD"
$bcde[47,2] = 246
$bcde[47,3] = 0.1629
$bcde[48,0] = "BRIGHOTEL"
$bcde[48,1] = "BRIGHOTEL"
$bcde[48,2] = 83.05
$bcde[48,3] = -0.1923
$bcde[49,0] = "EBGNG"
$bcde[49,1] = "EBGNG"
$bcde[49,2] = 337.9
$bcde[49,3] = 0.6254
$bcde[50,0] = "CHEMBONDCH"
$bcde[50,1] = "CHEMBONDCH"
$bcde[50,2] = 159.45
$bcde[50,3] = 2.2115
$bcde[51,0] = "ANTHEM"
$bcde[51,1] = "ANTHEM"
$bcde[51,2] = 703.6
$bcde[51,3] = 0.0711
$bcde[52,0] = "SMARTWORKS"
$bcde[52,1] = "SMARTWORKS"
$bcde[52,2] = 607.7
$bcde[52,3] = 2.2634
$bcde[53,0] = "TRAVELFOOD"
$bcde[53,1] = "TRAVELFOOD"
$bcde[53,2] = 1310
$bcde[53,3] = -0.3651
$bcde[54,0] = "CRIZAC"
$bcde[54,1] = "CRIZAC"
$bcde[54,2] = 307.5
$bcde[54,3] = -0.1947
$bcde[55,0] = "IGCL"
$bcde[55,1] = "IGCL"
$bcde[55,2] = 103.68
$bcde[55,3] = -0.3556
$bcde[56,0] = "SAMBHV"
$bcde[56,1] = "SAMBHV"
$bcde[56,2] = 114
$bcde[56,3] = -0.1227
$bcde[57,0] = "HDBFS"
$bcde[57,1] = "HDBFS"
$bcde[57,2] = 723.6
$bcde[57,3] = 0.0346
$bcde[58,0] = "ELLEN"
$bcde[58,1] = "ELLEN"
$bcde[58,2] = 464.65
$bcde[58,3] = 0
$bcde[59,0] = "KALPATARU"
$bcde[59,1] = "KALPATARU"
$bcde[59,2] = 384
$bcde[59,3] = -0.1171
$bcde[60,0] = "GLOBECIVIL"
$bcde[60,1] = "GLOBECIVIL"
$bcde[60,2] = 73.37
$bcde[60,3] = -0.9049
$bcde[61,0] = "RAYMONDREL"
$bcde[61,1] = "RAYMONDREL"
$bcde[61,2] = 620.3
$bcde[61,3] = -0.5053
$bcde[62,0] = "ARISINFRA"
$bcde[62,1] = "ARISINFRA"
$bcde[62,2] = 163.5
$bcde[62,3] = -0.0183
$bcde[63,0] = "ABLBL"
$bcde[63,1] = "ABLBL"
$bcde[63,2] = 136.89
$bcde[63,3] = 0.5066
$bcde[64,0] = "OSWALPUMPS"
$bcde[64,1] = "OSWALPUMPS"
$bcde[64,2] = 720
$bcde[64,3] = -0.1179
$bcde[65,0] = "ENRIN"
$bcde[65,1] = "ENRIN"
$bcde[65,2] = 3225
$bcde[65,3] = 0.5832
$bcde[66,0] = "BLUSPRING"
$bcde[66,1] = "BLUSPRING"
$bcde[66,2] = 78.64
$bcde[66,3] = 0.5112
$bcde[67,0] = "DIGITIDE"
$bcde[67,1] = "DIGITIDE"
$bcde[67,2] = 158.39
$bcde[67,3] = 0.5013
$bcde[68,0] = "SCODATUBES"
$bcde[68,1] = "SCODATUBES"
$bcde[68,2] = 160.25
$bcde[68,3] = 0.5017
$bcde[69,0] = "PROSTARM"
$bcde[69,1] = "PROSTARM"
$bcde[69,2] = 168.4
$bcde[69,3] = -0.1838
$bcde[70,0] = "THELEELA"
$bcde[70,1] = "THELEELA"
$bcde[70,2] = 439
$bcde[70,3] = 0.6996
$bcde[71,0] = "AEGISVOPAK"
$bcde[71,1] = "AEGISVOPAK"
$bcde[71,2] = 278.7
$bcde[71,3] = 0.5956
$bcde[72,0] = "BELRISE"
$bcde[72,1] = "BELRISE"
$bcde[72,2] = 145.95
$bcde[72,3] = 0.5789
$bcde[73,0] = "BORANA"
$bcde[73,1] = "BORANA"
$bcde[73,2] = 229.8
$bcde[73,3] = 0.5205
$bcde[74,0] = "ATHERENERG"
$bcde[74,1] = "ATHERENERG"
$bcde[74,2] = 691
$bcde[74,3] = -0.094

$ws.Range("B2:E76").Value = $bcde

$hcol = New-Object 'object[,]' 75,1
$hcol[0,0] = 11183.4
$hcol[1,0] = 10768.0802
$hcol[2,0] = 6181.9403
$hcol[3,0] = 111929.5674
$hcol[4,0] = 139571.2939
$hcol[5,0] = 8680.6865
$hcol[6,0] = 882.5144
$hcol[7,0] = 301.1595
$hcol[8,0] = 677.9608
$hcol[9,0] = 1051.7058
$hcol[10,0] = 4717.6466
$hcol[11,0] = 12861.3483
$hcol[12,0] = 2969.361
$hcol[13,0] = 600.5014
$hcol[14,0] = 6076.6178
$hcol[15,0] = 1397.7295
$hcol[16,0] = 2695.9529
$hcol[17,0] = 3218.119
$hcol[18,0] = 1140.8575
$hcol[19,0] = 7849.992
$hcol[20,0] = 4371.929
$hcol[21,0] = 6802.6598
$hcol[22,0] = 1511.9662
$hcol[23,0] = 350.7438
$hcol[24,0] = 3261.713
$hcol[25,0] = 2239.5386
$hcol[26,0] = 403.1382
$hcol[27,0] = 21827.1453
$hcol[28,0] = 469.1362
$hcol[29,0] = 564.254
$hcol[30,0] = 717.4921
$hcol[31,0] = 1302.7493
$hcol[32,0] = 2574.7248
$hcol[33,0] = 1283.5574
$hcol[34,0] = 4434.6349
$hcol[35,0] = 1156.8437
$hcol[36,0] = 740.3895
$hcol[37,0] = 11835.3813
$hcol[38,0] = 923.5866
$hcol[39,0] = 10759.6351
$hcol[40,0] = 1970.1464
$hcol[41,0] = 18491.3186
$hcol[42,0] = 0
$hcol[43,0] = 8658.718
$hcol[44,0] = 2592.8145
$hcol[45,0] = 770.3239
$hcol[46,0] = 15538.9329
$hcol[47,0] = 1770.6778
$hcol[48,0] = 3160.67
$hcol[49,0] = 3828.5067
$hcol[50,0] = 419.5866
$hcol[51,0] = 39486.8027
$hcol[52,0] = 6789.5693
$hcol[53,0] = 17313.2186
$hcol[54,0] = 5391.2108
$hcol[55,0] = 657.8457
$hcol[56,0] = 3363.3797
$hcol[57,0] = 60006.7092
$hcol[58,0] = 0
$hcol[59,0] = 7916.3635
$hcol[60,0] = 442.1595
$hcol[61,0] = 4150.5393
$hcol[62,0] = 1325.3865
$hcol[63,0] = 16620.4148
$hcol[64,0] = 8216.0619
$hcol[65,0] = 114182.9175
$hcol[66,0] = 1165.3802
$hcol[67,0] = 2347.4427
$hcol[68,0] = 955.2474
$hcol[69,0] = 993.2732
$hcol[70,0] = 14558.8937
$hcol[71,0] = 30696.9042
$hcol[72,0] = 12913.037
$hcol[73,0] = 609.1381
$hcol[74,0] = 26305.1261

$ws.Range("H2:H76").Value = $hcol

Write-Host "Edit complete"
